# "smaller window, fourth wall"
#
# The building's wall/roof/floor table gains a fourth wall (W4). In the
# original sheet the second data row (old row 5) represented a single
# "G1" window element whose width/height doubled as the gross opening for
# a wall that wasn't itself listed. This edit inserts a new row above it
# for that wall (W4, a Solid Wall w/In made of Concrete+Insulation), sizes
# its net surface as gross wall area minus the window's surface, and
# shrinks the window itself (G1) down a bit. Everything below ripples
# down one row (R1, R2, F1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 5; rows 5-8 (G1, R1, R2, F1) shift down to 6-9 and
# Excel auto-adjusts every relative formula/reference along the way
# (e.g. M2's "=F8*G8*H8" becomes "=F9*G9*H9", L6's own formula keeps its
# same-row relative refs, etc).
$ws.Rows("5:5").Insert()

# New row 5: the fourth wall, W4 - Solid Wall w/In (Concrete + Insulation).
# Its net surface subtracts out the window (now row 6) cut into it.
$ws.Range("A5").Value = "W4"
$ws.Range("B5").Value = "Solid Wall w/In"
$ws.Range("C5").Value = "Concrete"
$ws.Range("D5").Value = "Insulation"
$ws.Range("E5").Value = "N"
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = "N"
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = 0.1
$ws.Range("J5").Value = 0.1
$ws.Range("K5").Value = "N"
$ws.Range("L5").Formula = "=F5*H5 - L6"
$ws.Range("N5").Value = 90
$ws.Range("O5").Value = 90

# Row 6 (the G1 window, shifted down from the old row 5) shrinks from
# 3x3 to 2.5x2.5.
$ws.Range("F6").Value = 2.5
$ws.Range("H6").Value = 2.5

# Update the view's active cell/selection to match.
$ws.Range("H7").Select() | Out-Null
